$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.884.32"
$ws.Range("E2").Value = "  -0.85%  "
$ws.Range("D3").Value = "1.862.37"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.85"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3623"
$ws.Range("E8").Value = "  -3.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07169"
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8962"
$ws.Range("E10").Value = "  +1.12%  "
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.857.68"
$ws.Range("E12").Value = "  -0.43%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07471"
$ws.Range("E13").Value = "  -1.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.70"
$ws.Range("E14").Value = "  +3.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.230"
$ws.Range("E15").Value = "  -1.67%  "
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.16"
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").Value = "26.921.95"
$ws.Range("E20").Value = "  -0.88%  "
$ws.Range("E21").Value = "  -0.87%  "
$ws.Range("D22").Value = "2.083.91"
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.35"
$ws.Range("E23").Value = "  -1.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.416"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.80"
$ws.Range("E25").Value = "  -1.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.795"
$ws.Range("E26").Value = "  -2.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.86"
$ws.Range("E27").Value = "  -0.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.056"
$ws.Range("E28").Value = "  -1.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.88"
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.680"
$ws.Range("E30").Value = "  -1.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.678"
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09261"
$ws.Range("E32").Value = "  +2.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05093"
$ws.Range("E33").Value = "  -0.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.984"
$ws.Range("E34").Value = "  -3.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7428"
$ws.Range("E35").Value = "  +0.62%  "
$ws.Range("E36").Value = "  -0.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.278"
$ws.Range("E37").Value = "  +7.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02000"
$ws.Range("E38").Value = "  -1.60%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.502"
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5576"
$ws.Range("E40").Value = "  +4.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.072"
$ws.Range("E41").Value = "  -0.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "118.47"
$ws.Range("E42").Value = "  +2.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.483"
$ws.Range("E43").Value = "  -1.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.513"
$ws.Range("E44").Value = "  +2.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1468"
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4705"
$ws.Range("E46").Value = "  +1.63%  "
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.02"
$ws.Range("E48").Value = "  +0.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.565"
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.94"
$ws.Range("E50").Value = "  +1.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "62.97"
$ws.Range("E51").Value = "  -2.39%  "
